# Auto-generated edit script: update Goblin Market price/profit data per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 9).Value = 69.111115  # I5: was 69.888885
$ws.Cells.Item(5, 10).Value = 168  # J5: was 166.83333
$ws.Cells.Item(5, 11).Value = 69.111115  # K5: was 69.888885
$ws.Cells.Item(5, 12).Value = 168  # L5: was 166.83333
$ws.Cells.Item(5, 13).Value = 45.888885  # M5: was 45.111115
$ws.Cells.Item(5, 14).Value = -398  # N5: was -396.83333

$ws.Cells.Item(28, 8).Value = 9319.5  # H28: was 19152.637
$ws.Cells.Item(28, 9).Value = 15364.714  # I28: was 41090.4
$ws.Cells.Item(28, 10).Value = 856.2  # J28: was 871.1667
$ws.Cells.Item(28, 11).Value = 15364.714  # K28: was 41090.4
$ws.Cells.Item(28, 12).Value = 856.2  # L28: was 871.1667
$ws.Cells.Item(28, 13).Value = -14879.714  # M28: was -40605.4
$ws.Cells.Item(28, 14).Value = -1826.2  # N28: was -1841.1667

$ws.Cells.Item(41, 8).Value = 417  # H41: was 427.72726
$ws.Cells.Item(41, 9).Value = 138.16667  # I41: was 106
$ws.Cells.Item(41, 11).Value = 138.16667  # K41: was 106
$ws.Cells.Item(41, 13).Value = 301.83333  # M41: was 334

$ws.Cells.Item(43, 8).Value = 8912.333000000001  # H43: was 5805.231
$ws.Cells.Item(43, 9).Value = 3750  # I43: was 3098.5
$ws.Cells.Item(43, 10).Value = 11493.5  # J43: was 7008.222
$ws.Cells.Item(43, 11).Value = 3750  # K43: was 3098.5
$ws.Cells.Item(43, 12).Value = 11493.5  # L43: was 7008.222
$ws.Cells.Item(43, 13).Value = -3681  # M43: was -3029.5
$ws.Cells.Item(43, 14).Value = -11631.5  # N43: was -7146.222

$ws.Cells.Item(88, 8).Value = 4024.7856  # H88: was 3719.2354
$ws.Cells.Item(88, 9).Value = 0  # I88: was 800
$ws.Cells.Item(88, 10).Value = 4024.7856  # J88: was 3901.6875
$ws.Cells.Item(88, 11).Value = 0  # K88: was 800
$ws.Cells.Item(88, 12).Value = 4024.7856  # L88: was 3901.6875
$ws.Cells.Item(88, 13).ClearContents()  # M88: was -394
$ws.Cells.Item(88, 14).Value = -4836.7856  # N88: was -4713.6875

$ws.Cells.Item(91, 8).Value = 4024.7856  # H91: was 3719.2354
$ws.Cells.Item(91, 9).Value = 0  # I91: was 800
$ws.Cells.Item(91, 10).Value = 4024.7856  # J91: was 3901.6875
$ws.Cells.Item(91, 11).Value = 0  # K91: was 800
$ws.Cells.Item(91, 12).Value = 4024.7856  # L91: was 3901.6875
$ws.Cells.Item(91, 13).ClearContents()  # M91: was 604
$ws.Cells.Item(91, 14).Value = -6832.7856  # N91: was -6709.6875

$ws.Cells.Item(92, 8).Value = 2674  # H92: was 3055.8333
$ws.Cells.Item(92, 9).Value = 2485.8572  # I92: was 3000.4
$ws.Cells.Item(92, 10).Value = 3332.5  # J92: was 3333
$ws.Cells.Item(92, 11).Value = 2485.8572  # K92: was 3000.4
$ws.Cells.Item(92, 12).Value = 3332.5  # L92: was 3333
$ws.Cells.Item(92, 13).Value = -1237.8572  # M92: was -1752.4
$ws.Cells.Item(92, 14).Value = -5828.5  # N92: was -5829

$ws.Cells.Item(93, 8).Value = 59152.75  # H93: was 65537
$ws.Cells.Item(93, 10).Value = 59152.75  # J93: was 65537
$ws.Cells.Item(93, 12).Value = 59152.75  # L93: was 65537
$ws.Cells.Item(93, 14).Value = -64144.75  # N93: was -70529

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 0  # H3: was 5000
$ws.Cells.Item(3, 9).Value = 0  # I3: was 5000
$ws.Cells.Item(3, 11).Value = 0  # K3: was 5000
$ws.Cells.Item(3, 13).ClearContents()  # M3: was -4885

$ws.Cells.Item(54, 8).Value = 20000  # H54: was 0
$ws.Cells.Item(54, 10).Value = 20000  # J54: was 0
$ws.Cells.Item(54, 12).Value = 20000  # L54: was 0
$ws.Cells.Item(54, 14).Value = -21538  # N54: was None

$ws.Cells.Item(74, 8).Value = 2030.1333  # H74: was 1663.8096
$ws.Cells.Item(74, 9).Value = 2068  # I74: was 1681.1052
$ws.Cells.Item(74, 10).Value = 1500  # J74: was 1499.5
$ws.Cells.Item(74, 11).Value = 2068  # K74: was 1681.1052
$ws.Cells.Item(74, 12).Value = 1500  # L74: was 1499.5
$ws.Cells.Item(74, 13).Value = -1194  # M74: was -807.1052
$ws.Cells.Item(74, 14).Value = -3248  # N74: was -3247.5

$ws.Cells.Item(77, 8).Value = 2030.1333  # H77: was 1663.8096
$ws.Cells.Item(77, 9).Value = 2068  # I77: was 1681.1052
$ws.Cells.Item(77, 10).Value = 1500  # J77: was 1499.5
$ws.Cells.Item(77, 11).Value = 10340  # K77: was 8405.526
$ws.Cells.Item(77, 12).Value = 7500  # L77: was 7497.5
$ws.Cells.Item(77, 13).Value = -5972  # M77: was -4037.526
$ws.Cells.Item(77, 14).Value = -16236  # N77: was -16233.5

$ws.Cells.Item(94, 8).Value = 20299.166  # H94: was 24959
$ws.Cells.Item(94, 10).Value = 20299.166  # J94: was 24959
$ws.Cells.Item(94, 12).Value = 20299.166  # L94: was 24959
$ws.Cells.Item(94, 14).Value = -22101.166  # N94: was -26761

$ws.Cells.Item(97, 8).Value = 299.63635  # H97: was 323.95456
$ws.Cells.Item(97, 9).Value = 306.5263  # I97: was 334.6842
$ws.Cells.Item(97, 11).Value = 306.5263  # K97: was 334.6842
$ws.Cells.Item(97, 13).Value = 189.4737  # M97: was 161.3158

$ws.Cells.Item(102, 8).Value = 3040.3235  # H102: was 3313.24
$ws.Cells.Item(102, 9).Value = 1840.4138  # I102: was 1641.6
$ws.Cells.Item(102, 11).Value = 1840.4138  # K102: was 1641.6
$ws.Cells.Item(102, 13).Value = -218.4138  # M102: was -19.59999999999991

$ws.Cells.Item(104, 8).Value = 5080.6  # H104: was 5640.6
$ws.Cells.Item(104, 10).Value = 5599.5  # J104: was 6999.5
$ws.Cells.Item(104, 12).Value = 5599.5  # L104: was 6999.5
$ws.Cells.Item(104, 14).Value = -12587.5  # N104: was -13987.5

$ws.Cells.Item(110, 8).Value = 1006.1923  # H110: was 1016.44
$ws.Cells.Item(110, 10).Value = 1085.3334  # J110: was 1152.4
$ws.Cells.Item(110, 12).Value = 1085.3334  # L110: was 1152.4
$ws.Cells.Item(110, 14).Value = -5175.3334  # N110: was -5242.4

$ws.Cells.Item(122, 8).Value = 22225714  # H122: was 9261948
$ws.Cells.Item(122, 9).Value = 111111110  # I122: was 13890740
$ws.Cells.Item(122, 11).Value = 333333330  # K122: was 41672220
$ws.Cells.Item(122, 13).Value = -333330880  # M122: was -41669770

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(10, 8).Value = 650.8333  # H10: was 713.75
$ws.Cells.Item(10, 9).Value = 101.666664  # I10: was 127.5
$ws.Cells.Item(10, 10).Value = 1200  # J10: was 1300
$ws.Cells.Item(10, 11).Value = 101.666664  # K10: was 127.5
$ws.Cells.Item(10, 12).Value = 1200  # L10: was 1300
$ws.Cells.Item(10, 13).Value = 38.333336  # M10: was 12.5
$ws.Cells.Item(10, 14).Value = -1480  # N10: was -1580

$ws.Cells.Item(105, 8).Value = 1644.2941  # H105: was 1748.9131
$ws.Cells.Item(105, 9).Value = 1661.8372  # I105: was 1790.7894
$ws.Cells.Item(105, 11).Value = 1661.8372  # K105: was 1790.7894
$ws.Cells.Item(105, 13).Value = 85.16280000000006  # M105: was -43.78939999999989

$ws.Cells.Item(134, 8).Value = 4559.25  # H134: was 4052
$ws.Cells.Item(134, 9).Value = 4559.25  # I134: was 4052
$ws.Cells.Item(134, 11).Value = 13677.75  # K134: was 12156
$ws.Cells.Item(134, 13).Value = -11142.75  # M134: was -9621

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 107.72727  # H7: was 130.65
$ws.Cells.Item(7, 9).Value = 53.083332  # I7: was 57.81818
$ws.Cells.Item(7, 10).Value = 173.3  # J7: was 219.66667
$ws.Cells.Item(7, 11).Value = 53.083332  # K7: was 57.81818
$ws.Cells.Item(7, 12).Value = 173.3  # L7: was 219.66667
$ws.Cells.Item(7, 13).Value = 59.916668  # M7: was 55.18182
$ws.Cells.Item(7, 14).Value = -399.3  # N7: was -445.66667

$ws.Cells.Item(22, 8).Value = 1426.5  # H22: was 2044
$ws.Cells.Item(22, 9).Value = 752.5714  # I22: was 1089
$ws.Cells.Item(22, 11).Value = 752.5714  # K22: was 1089
$ws.Cells.Item(22, 13).Value = -402.5714  # M22: was -739

$ws.Cells.Item(38, 8).Value = 4956  # H38: was 3902
$ws.Cells.Item(38, 10).Value = 8771  # J38: was 6136
$ws.Cells.Item(38, 12).Value = 8771  # L38: was 6136
$ws.Cells.Item(38, 14).Value = -9525  # N38: was -6890

$ws.Cells.Item(42, 8).Value = 7770.8  # H42: was 8109.1665
$ws.Cells.Item(42, 10).Value = 9899  # J42: was 9866.333000000001
$ws.Cells.Item(42, 12).Value = 9899  # L42: was 9866.333000000001
$ws.Cells.Item(42, 14).Value = -11085  # N42: was -11052.333

$ws.Cells.Item(46, 8).Value = 4956  # H46: was 3902
$ws.Cells.Item(46, 10).Value = 8771  # J46: was 6136
$ws.Cells.Item(46, 12).Value = 8771  # L46: was 6136
$ws.Cells.Item(46, 14).Value = -9193  # N46: was -6558

$ws.Cells.Item(107, 8).Value = 1263.7317  # H107: was 1268.9269
$ws.Cells.Item(107, 9).Value = 455.07693  # I107: was 473.35715
$ws.Cells.Item(107, 10).Value = 1639.1786  # J107: was 1681.4445
$ws.Cells.Item(107, 11).Value = 455.07693  # K107: was 473.35715
$ws.Cells.Item(107, 12).Value = 1639.1786  # L107: was 1681.4445
$ws.Cells.Item(107, 13).Value = 1464.92307  # M107: was 1446.64285
$ws.Cells.Item(107, 14).Value = -5479.1786  # N107: was -5521.4445

$ws.Cells.Item(120, 8).Value = 39500  # H120: was 40000
$ws.Cells.Item(120, 10).Value = 39500  # J120: was 40000
$ws.Cells.Item(120, 12).Value = 39500  # L120: was 40000
$ws.Cells.Item(120, 14).Value = -46758  # N120: was -47258

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(104, 8).Value = 21217.812  # H104: was 14767.923
$ws.Cells.Item(104, 9).Value = 4000  # I104: was 3499
$ws.Cells.Item(104, 10).Value = 22365.666  # J104: was 16816.818
$ws.Cells.Item(104, 11).Value = 12000  # K104: was 10497
$ws.Cells.Item(104, 12).Value = 67096.99800000001  # L104: was 50450.454
$ws.Cells.Item(104, 13).Value = -9379  # M104: was -7876
$ws.Cells.Item(104, 14).Value = -72338.99800000001  # N104: was -55692.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 68.40000000000001  # H2: was 64.25
$ws.Cells.Item(2, 10).Value = 73.333336  # J2: was 50
$ws.Cells.Item(2, 12).Value = 73.333336  # L2: was 50
$ws.Cells.Item(2, 14).Value = -299.333336  # N2: was -276

$ws.Cells.Item(24, 8).Value = 28333.334  # H24: was 28750
$ws.Cells.Item(24, 10).Value = 28333.334  # J24: was 28750
$ws.Cells.Item(24, 12).Value = 28333.334  # L24: was 28750
$ws.Cells.Item(24, 14).Value = -28679.334  # N24: was -29096

$ws.Cells.Item(70, 8).Value = 22225936  # H70: was 20837802
$ws.Cells.Item(70, 9).Value = 30306310  # I70: was 37041380
$ws.Cells.Item(70, 10).Value = 4908.25  # J70: was 4631.857
$ws.Cells.Item(70, 11).Value = 30306310  # K70: was 37041380
$ws.Cells.Item(70, 12).Value = 4908.25  # L70: was 4631.857
$ws.Cells.Item(70, 13).Value = -30306040  # M70: was -37041110
$ws.Cells.Item(70, 14).Value = -5448.25  # N70: was -5171.857

$ws.Cells.Item(73, 8).Value = 22225936  # H73: was 20837802
$ws.Cells.Item(73, 9).Value = 30306310  # I73: was 37041380
$ws.Cells.Item(73, 10).Value = 4908.25  # J73: was 4631.857
$ws.Cells.Item(73, 11).Value = 30306310  # K73: was 37041380
$ws.Cells.Item(73, 12).Value = 4908.25  # L73: was 4631.857
$ws.Cells.Item(73, 13).Value = -30305374  # M73: was -37040444
$ws.Cells.Item(73, 14).Value = -6780.25  # N73: was -6503.857

$ws.Cells.Item(97, 8).Value = 392  # H97: was 279.42307
$ws.Cells.Item(97, 9).Value = 335.58334  # I97: was 286.08334
$ws.Cells.Item(97, 10).Value = 1069  # J97: was 199.5
$ws.Cells.Item(97, 11).Value = 335.58334  # K97: was 286.08334
$ws.Cells.Item(97, 12).Value = 1069  # L97: was 199.5
$ws.Cells.Item(97, 13).Value = 160.41666  # M97: was 209.91666
$ws.Cells.Item(97, 14).Value = -2061  # N97: was -1191.5

$ws.Cells.Item(107, 8).Value = 576.1818  # H107: was 552.5333000000001
$ws.Cells.Item(107, 10).Value = 1445.5  # J107: was 806.8333
$ws.Cells.Item(107, 12).Value = 1445.5  # L107: was 806.8333
$ws.Cells.Item(107, 14).Value = -5285.5  # N107: was -4646.8333

$ws.Cells.Item(122, 8).Value = 9073.806  # H122: was 8828.621999999999
$ws.Cells.Item(122, 9).Value = 9689.571  # I122: was 9143.634
$ws.Cells.Item(122, 10).Value = 6918.625  # J122: was 7478.5713
$ws.Cells.Item(122, 11).Value = 29068.713  # K122: was 27430.902
$ws.Cells.Item(122, 12).Value = 20755.875  # L122: was 22435.7139
$ws.Cells.Item(122, 13).Value = -26618.713  # M122: was -24980.902
$ws.Cells.Item(122, 14).Value = -25655.875  # N122: was -27335.7139

$ws.Cells.Item(123, 8).Value = 79999  # H123: was 39999
$ws.Cells.Item(123, 10).Value = 79999  # J123: was 39999
$ws.Cells.Item(123, 12).Value = 79999  # L123: was 39999
$ws.Cells.Item(123, 14).Value = -84899  # N123: was -44899

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6416.6665  # H7: was 5875
$ws.Cells.Item(7, 9).Value = 5700  # I7: was 5875
$ws.Cells.Item(7, 10).Value = 10000  # J7: was 0
$ws.Cells.Item(7, 11).Value = 5700  # K7: was 5875
$ws.Cells.Item(7, 12).Value = 10000  # L7: was 0
$ws.Cells.Item(7, 13).Value = -5588  # M7: was -5763
$ws.Cells.Item(7, 14).Value = -10224  # N7: was None

$ws.Cells.Item(22, 8).Value = 3250.9666  # H22: was 3401.0715
$ws.Cells.Item(22, 10).Value = 3515.7778  # J22: was 3811.5625
$ws.Cells.Item(22, 12).Value = 3515.7778  # L22: was 3811.5625
$ws.Cells.Item(22, 14).Value = -4105.7778  # N22: was -4401.5625

$ws.Cells.Item(27, 8).Value = 3250.9666  # H27: was 3401.0715
$ws.Cells.Item(27, 10).Value = 3515.7778  # J27: was 3811.5625
$ws.Cells.Item(27, 12).Value = 3515.7778  # L27: was 3811.5625
$ws.Cells.Item(27, 14).Value = -3729.7778  # N27: was -4025.5625

$ws.Cells.Item(61, 8).Value = 6039.7334  # H61: was 6299.2856
$ws.Cells.Item(61, 9).Value = 4732.8887  # I61: was 5023.75
$ws.Cells.Item(61, 11).Value = 4732.8887  # K61: was 5023.75
$ws.Cells.Item(61, 13).Value = -4530.8887  # M61: was -4821.75

$ws.Cells.Item(93, 8).Value = 5334.75  # H93: was 4833.5
$ws.Cells.Item(93, 9).Value = 782.6667  # I93: was 799
$ws.Cells.Item(93, 10).Value = 6385.231  # J93: was 5640.4
$ws.Cells.Item(93, 11).Value = 782.6667  # K93: was 799
$ws.Cells.Item(93, 12).Value = 6385.231  # L93: was 5640.4
$ws.Cells.Item(93, 13).Value = 465.3333  # M93: was 449
$ws.Cells.Item(93, 14).Value = -8881.231  # N93: was -8136.4

$ws.Cells.Item(113, 8).Value = 6039.7334  # H113: was 6299.2856
$ws.Cells.Item(113, 9).Value = 4732.8887  # I113: was 5023.75
$ws.Cells.Item(113, 11).Value = 4732.8887  # K113: was 5023.75
$ws.Cells.Item(113, 13).Value = -2562.8887  # M113: was -2853.75

$ws.Cells.Item(126, 8).Value = 6416.6665  # H126: was 5875
$ws.Cells.Item(126, 9).Value = 5700  # I126: was 5875
$ws.Cells.Item(126, 10).Value = 10000  # J126: was 0
$ws.Cells.Item(126, 11).Value = 17100  # K126: was 17625
$ws.Cells.Item(126, 12).Value = 30000  # L126: was 0
$ws.Cells.Item(126, 13).Value = -14630  # M126: was -15155
$ws.Cells.Item(126, 14).Value = -34940  # N126: was None

$ws.Cells.Item(132, 8).Value = 5525.2104  # H132: was 5733
$ws.Cells.Item(132, 9).Value = 5310.875  # I132: was 5699
$ws.Cells.Item(132, 10).Value = 6668.3335  # J132: was 6005
$ws.Cells.Item(132, 11).Value = 15932.625  # K132: was 17097
$ws.Cells.Item(132, 12).Value = 20005.0005  # L132: was 18015
$ws.Cells.Item(132, 13).Value = -13402.625  # M132: was -14567
$ws.Cells.Item(132, 14).Value = -25065.0005  # N132: was -23075

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 27385.666  # H95: was 29867.8
$ws.Cells.Item(95, 10).Value = 27385.666  # J95: was 29867.8
$ws.Cells.Item(95, 12).Value = 27385.666  # L95: was 29867.8
$ws.Cells.Item(95, 14).Value = -32877.666  # N95: was -35359.8

$ws.Cells.Item(107, 8).Value = 1466.2222  # H107: was 2140
$ws.Cells.Item(107, 9).Value = 599.25  # I107: was 600
$ws.Cells.Item(107, 10).Value = 2159.8  # J107: was 2525
$ws.Cells.Item(107, 11).Value = 1797.75  # K107: was 1800
$ws.Cells.Item(107, 12).Value = 6479.400000000001  # L107: was 7575
$ws.Cells.Item(107, 13).Value = 122.25  # M107: was 120
$ws.Cells.Item(107, 14).Value = -10319.4  # N107: was -11415

$ws.Cells.Item(132, 8).Value = 4537.3145  # H132: was 4549.457
$ws.Cells.Item(132, 9).Value = 4736.1  # I132: was 4647.839
$ws.Cells.Item(132, 10).Value = 4736.1  # J132: was 3787
$ws.Cells.Item(132, 11).Value = 14208.3  # K132: was 13943.517
$ws.Cells.Item(132, 12).Value = 10033.8  # L132: was 11361
$ws.Cells.Item(132, 13).Value = -11678.3  # M132: was -11413.517
$ws.Cells.Item(132, 14).Value = -15093.8  # N132: was -16421
